# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" counts (column F) and two refreshed
# cover-image URLs (column I) on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 1892
$ws1.Range("F6").Value  = 2656
$ws1.Range("F7").Value  = 181
$ws1.Range("F9").Value  = 179
$ws1.Range("F10").Value = 1559
$ws1.Range("F11").Value = 543
$ws1.Range("F15").Value = 24
$ws1.Range("F17").Value = 215
$ws1.Range("F18").Value = 222
$ws1.Range("F21").Value = 197
$ws1.Range("F23").Value = 1701
$ws1.Range("F25").Value = 416
$ws1.Range("F26").Value = 34
$ws1.Range("F30").Value = 433
$ws1.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg"

# Sheet 4: 全部类型 (all types) - same underlying events, rows shifted by one
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 1892
$ws4.Range("F5").Value  = 169
$ws4.Range("F7").Value  = 2656
$ws4.Range("F8").Value  = 181
$ws4.Range("F10").Value = 179
$ws4.Range("F11").Value = 1559
$ws4.Range("F12").Value = 543
$ws4.Range("F16").Value = 24
$ws4.Range("F18").Value = 215
$ws4.Range("F19").Value = 222
$ws4.Range("F22").Value = 197
$ws4.Range("F24").Value = 1701
$ws4.Range("F26").Value = 416
$ws4.Range("F27").Value = 34
$ws4.Range("F31").Value = 433
$ws4.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg"
